$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs","Sema7a","Plxnc1","ECs",3,1,23.128865,69.386595,0.7917836846260858,0.7917836846260858,2,0.6666666666666666,23.76805833333333,71.304175,0.8847414975178436,0.8847414975178437,549.7282125037917,4947.553912534125,0.7005238828462792,0.7005238828462793),
    @("ECs","Sema7a","Plxnc1","FAPs",3,1,23.128865,69.386595,0.7917836846260858,0.7917836846260858,3,1,2.050040666666666,6.150122,0.07631065289230869,0.0763106528923087,47.41511382384333,426.73602441459,0.06042152992329445,0.06042152992329446),
    @("ECs","Sema7a","Plxnc1","sCs",3,1,23.128865,69.386595,0.7917836846260858,0.7917836846260858,3,1,1.046311,3.138933,0.03894784958984768,0.03894784958984768,24.199985867015,217.799872803135,0.03083827185651218,0.03083827185651218),
    @("FAPs","Sema7a","Plxnc1","ECs",3,1,1.158105666666667,3.474317,0.03964609469334889,0.03964609469334889,2,0.6666666666666666,23.76805833333333,71.304175,0.8847414975178436,0.8847414975178437,27.52592304149722,247.733307373475,0.03507654518972773,0.03507654518972773),
    @("FAPs","Sema7a","Plxnc1","FAPs",3,1,1.158105666666667,3.474317,0.03964609469334889,0.03964609469334889,3,1,2.050040666666666,6.150122,0.07631065289230869,0.0763106528923087,2.374163712963778,21.367473416674,0.003025419370679749,0.00302541937067975),
    @("FAPs","Sema7a","Plxnc1","sCs",3,1,1.158105666666667,3.474317,0.03964609469334889,0.03964609469334889,3,1,1.046311,3.138933,0.03894784958984768,0.03894784958984768,1.211738698195667,10.905648283761,0.001544130132941411,0.001544130132941411),
    @("sCs","Sema7a","Plxnc1","ECs",3,1,4.924119999999999,14.77236,0.1685702206805652,0.1685702206805652,2,0.6666666666666666,23.76805833333333,71.304175,0.8847414975178436,0.8847414975178437,117.0367714003333,1053.330942603,0.1491410694818366,0.1491410694818367),
    @("sCs","Sema7a","Plxnc1","FAPs",3,1,4.924119999999999,14.77236,0.1685702206805652,0.1685702206805652,3,1,2.050040666666666,6.150122,0.07631065289230869,0.0763106528923087,10.09464624754666,90.85181622791998,0.01286370359833449,0.01286370359833449),
    @("sCs","Sema7a","Plxnc1","sCs",3,1,4.924119999999999,14.77236,0.1685702206805652,0.1685702206805652,3,1,1.046311,3.138933,0.03894784958984768,0.03894784958984768,5.152160921319999,46.36944829187999,0.006565447600394086,0.006565447600394086)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
